# modificacion para tesoreria y comercializacion
# Adds a new payment record (row 30) to the Worksheet, mirroring the
# existing columns: FECHA, RUBRO, UNIDAD, Nº RECIBO, NOMBRE Y APELLIDOS,
# DETALLE DE PAGO, MONTO, USUARIO SISTEMA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (style/borders) from the previous data row (29) so the
# new row (30) matches the rest of the table instead of picking up a
# default/blank style.
$srcRange = $ws.Range("A29:H29")
$dstRange = $ws.Range("A30:H30")
$srcRange.Copy($dstRange)

$ws.Cells.Item(30, 1).Value = "18/05/2021"
$ws.Cells.Item(30, 2).Value = 12200
$ws.Cells.Item(30, 3).Value = "Ing. Ecopiscicultura"
$ws.Cells.Item(30, 4).Value = 17
$ws.Cells.Item(30, 5).Value = "ACOSTA SAAVEDRA MARIA SOLEDAD"
$ws.Cells.Item(30, 6).Value = "Pago de Alevines"
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 8).Value = "finley_1231@hotmail.com"

# The longer name in column E forces the column to widen to fit it,
# same as Excel's bestFit/AutoFit recompute on the new content.
$ws.Columns.Item(5).ColumnWidth = 34.5

# Keep the whole-table selection in sync with the grown used range.
$ws.Range("A1:H30").Select() | Out-Null
